# Applies updated calibration values to the workbook.
$wb = $excel.ActiveWorkbook

# --- start_price sheet ---
$wsStartPrice = $wb.Worksheets.Item("start_price")
$wsStartPrice.Range("A2").Value = 43.54

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.00388368359661166
$wsLinear.Range("B3").Value = -0.07014327697295385
$wsLinear.Range("B4").Value = 0.1765392601810676
$wsLinear.Range("B5").Value = "[1.0, 0.1700418353774329, 0.043085070444321506, 0.017879366349800793, 0.028310860240480212, 0.0051177370905848155, 0.0636560262612801, 0.12983819500731605, 0.03982696699992425, -0.0033214840383491654, -0.0468891723682396, -0.04693992374192806, -0.06216728036303092, 0.007607422620499608, 0.06270798521092576, 0.017838159087306172, -0.04450874654420988, -0.0564587193032945, -0.005173234393754638, -0.032868203981919177]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B3").Value = 0.7989210397253556
$wsNonLinear.Range("B4").Value = -0.006874752269550893
$wsNonLinear.Range("B5").Value = -0.1301792071301368
$wsNonLinear.Range("B6").Value = 0.1909290066002897
$wsNonLinear.Range("B7").Value = 0.003784768292128477
$wsNonLinear.Range("B8").Value = -0.04472842887261273
$wsNonLinear.Range("B9").Value = 0.1651283607289172
$wsNonLinear.Range("B10").Value = "[1.0, 0.17181173358083407, 0.044129666261103005, 0.01733360681727324, 0.028525435878691024, 0.005521722549867253, 0.06315522974829448, 0.1303800381036637, 0.039975987744650845, -0.0036819442745926515, -0.048166861141688, -0.04842945972245285, -0.06267442031223332, 0.007246887921281102, 0.06316503076489494, 0.018118055577219405, -0.04524153417354356, -0.05764182531529819, -0.006389029750875455, -0.033396727318463125]"
